$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.828.54'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '1.875.93'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').Value = '0.9976'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '242.74'
$ws.Range('E5').Value = '  -3.08%  '
$ws.Range('D6').Value = '0.9969'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').Value = '0.4924'
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('D8').Value = '44.25'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').Value = '0.2890'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = '0.06587'
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').Value = '1.874.85'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '16.85'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').Value = '0.07167'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '0.6667'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '85.42'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '4.795'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').Value = '29.838.23'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '0.000007824'
$ws.Range('E18').Value = '  +4.76%  '
$ws.Range('D19').Value = '0.9966'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '12.73'
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('D21').Value = '2.117.17'
$ws.Range('E21').Value = '  +1.89%  '
$ws.Range('D22').Value = '0.9964'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '4.735'
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '5.560'
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '9.084'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('D26').Value = '146.80'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('D27').Value = '133.83'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').Value = '16.68'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').Value = '1.921'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').Value = '1.375'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').Value = '4.160'
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('D32').Value = '0.08572'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').Value = '3.908'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('D34').Value = '0.04977'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.107'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7033'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').Value = '2.661'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').Value = '2.197'
$ws.Range('E38').Value = '  -5.64%  '
$ws.Range('D39').Value = '2.665'
$ws.Range('E39').Value = '  -2.77%  '
$ws.Range('D40').Value = '0.9293'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('D41').Value = '0.01633'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '6.048'
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('D43').Value = '0.9934'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').Value = '102.49'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').Value = '0.4162'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '7.519'
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('D47').Value = '0.1255'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('D49').Value = '32.46'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').Value = '8.175'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '0.3694'
$ws.Range('E51').Value = '  -0.13%  '
